$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3: Date 2012-10-28, Anlass "GD", Ort "Kirche", Kirchgemeindetag "8, 88, 182"
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A3").Value = 41210

$ws.Range("B3").Value = "GD"
$ws.Range("C3").Value = "Kirche"
$ws.Range("D3").Value = "8, 88, 182"

$ws.Range("D3").Select()
